$wb = $excel.ActiveWorkbook

# --- 1) Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is referenced by Overview!E2, Overview!F2, zh-cn!C2 and de-de!C2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the "zh-cn" / "de-de" / "Status" columns ---
# Original raw width 17.2159881591797 -> new raw width 13.4101845877511
# which corresponds to a COM ColumnWidth of 12.5 (closest achievable value)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
